$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: was J3 / Conn_01x04_Male -> now J1 with full connector data ---
# (H8 is plain text here, not a live hyperlink - matches target which has no style/relationship for it)
$ws.Range("A8").Value2 = "J1"
$ws.Range("B8").Value2 = "B3B-XH-A(LF)(SN)"
$ws.Range("H8").Value2 = "https://www.digikey.com/en/products/detail/jst-sales-america-inc/B3B-XH-A-LF-SN/1651046"
$ws.Range("C8").Value2 = "CONN HEADER VERT 3POS 2.5MM"
$ws.Range("D8").Value2 = "Digikey"
$ws.Range("E8").Value2 = 0.19
$ws.Range("F8").Value2 = 1

# --- Row 13: was J1 / Conn_01x03 -> now J3 with full connector data ---
$ws.Range("A13").Value2 = "J3"
$ws.Range("B13").Value2 = "B4B-XH-A(LF)(SN)"
$ws.Range("C13").Value2 = "CONN HEADER VERT 4POS 2.5MM"
$ws.Range("D13").Value2 = "Digikey"
$ws.Range("E13").Value2 = 0.21
$ws.Range("F13").Value2 = 1

# --- Row 15: was U2 / STDC14 -> now U2 with full connector data ---
$ws.Range("H15").Value2 = "https://www.digikey.com/en/products/detail/samtec-inc/FTSH-107-01-L-DV-K/6678186"
$ws.Range("C15").Value2 = "CONN HEADER SMD 14POS 1.27MM"
$ws.Range("B15").Value2 = "FTSH-107-01-L-DV-K"
$ws.Range("B15").WrapText = $true
$ws.Range("C15").WrapText = $true
$ws.Range("D15").Value2 = "Digikey"
$ws.Range("E15").Value2 = 5.72
$ws.Range("F15").Value2 = 1
$ws.Hyperlinks.Add($ws.Range("H15"), "https://www.digikey.com/en/products/detail/samtec-inc/FTSH-107-01-L-DV-K/6678186", [Type]::Missing, [Type]::Missing, "https://www.digikey.com/en/products/detail/samtec-inc/FTSH-107-01-L-DV-K/6678186")
$ws.Range("H15").Style = "Hyperlink"

# --- Column C width: widen to fit the longer descriptions now present ---
$ws.Columns.Item(3).ColumnWidth = 70.35

# --- Selection moved to C17 per the saved view state ---
$ws.Range("C17").Select()
